# "AUTO USATE" schema sheet: add the column catalogue (name / type / PK) for
# the used-cars table, in columns C..S, mirroring the existing B-column
# labels (nome colonna / tipo di dato / chiave primaria / attributi/proprietà).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered left-to-right column definitions: address, column name, sql type, is primary key
$columns = @(
    @{ Col = "C"; Name = "ID";                   Type = "SMALLINT";     Pk = $true  },
    @{ Col = "D"; Name = "marca";                 Type = "VARCHAR(50)";  Pk = $false },
    @{ Col = "E"; Name = "modello";                Type = "VARCHAR(100)"; Pk = $false },
    @{ Col = "F"; Name = "anno_immatricolazione";  Type = "DATE";         Pk = $false },
    @{ Col = "G"; Name = "km_percorsi";            Type = "INT";          Pk = $false },
    @{ Col = "H"; Name = "prezzo";                 Type = "DECIMAL(9,2)"; Pk = $false },
    @{ Col = "I"; Name = "tipo_auto";              Type = "VARCHAR(50)";  Pk = $false },
    @{ Col = "J"; Name = "tipo_cambio";            Type = "VARCHAR(50)";  Pk = $false },
    @{ Col = "K"; Name = "colore";                 Type = "VARCHAR(50)";  Pk = $false },
    @{ Col = "L"; Name = "potenza_motore";         Type = "SMALLINT";     Pk = $false },
    @{ Col = "M"; Name = "grandezza_serbatoio";    Type = "SMALLINT";     Pk = $false },
    @{ Col = "N"; Name = "numero_porte";           Type = "TINYINT";      Pk = $false },
    @{ Col = "O"; Name = "chiusura_elettrica";     Type = "BOOL";         Pk = $false },
    @{ Col = "P"; Name = "finestrini_elettrici";   Type = "BOOL";         Pk = $false },
    @{ Col = "Q"; Name = "ruota_di_scorta";        Type = "BOOL";         Pk = $false },
    @{ Col = "R"; Name = "radio";                  Type = "BOOL";         Pk = $false },
    @{ Col = "S"; Name = "bluetooth";              Type = "BOOL";         Pk = $false }
)

# Column widths tuned so the host's pixel-quantized ColumnWidth model lands as
# close as possible to the authored (Excel bestFit / manual) stored width.
$widths = @{
    "C" = 8.666666666666666
    "D" = 12.0
    "E" = 13.0
    "F" = 21.166666666666668
    "G" = 11.333333333333334
    "H" = 12.666666666666666
    "I" = 11.833333333333334
    "J" = 12.333333333333334
    "K" = 12.0
    "L" = 14.666666666666666
    "M" = 18.166666666666668
    "N" = 12.666666666666666
    "O" = 15.666666666666666
    "P" = 16.0
    "Q" = 13.666666666666666
    "S" = 9.0
}

foreach ($c in $columns) {
    # Row 2: column name, bold (header)
    $cell = $ws.Range($c.Col + "2")
    $cell.Value = $c.Name
    $cell.Font.Bold = $true

    # Row 3: data type
    $ws.Range($c.Col + "3").Value = $c.Type

    # Row 4: primary key marker
    if ($c.Pk) {
        $ws.Range($c.Col + "4").Value = "X"
    }

    # Column width (column R / "radio" keeps the sheet default — no entry in $widths)
    if ($widths.ContainsKey($c.Col)) {
        $ws.Range($c.Col + "1").ColumnWidth = $widths[$c.Col]
    }
}

# Selection left by the author after finishing the edit
$ws.Range("Q4").Select()

# Page setup (matches paperSize=9 / portrait added to the sheet)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
